$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 403, pushing existing rows 403:475 down to 404:476
$ws.Rows(403).Insert()

# Populate the new row 403 with the latest weekly price record.
# Columns A,B,C,E,F,G,H,I,R stay constant across this block of records.
$ws.Range("A403").Value = 4
$ws.Range("B403").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C403").Value = "Los Lagos"
$ws.Range("D403").Value = 45211
$ws.Range("E403").Value = 10
$ws.Range("F403").Value = 100112032
$ws.Range("G403").Value = "Zapallo italiano"
$ws.Range("H403").Value = "Sin especificar"
$ws.Range("I403").Value = "Primera"
$ws.Range("J403").Value = 150
$ws.Range("K403").Value = 26000
$ws.Range("L403").Value = 26000
$ws.Range("M403").Value = 26000
$ws.Range("N403").Value = "$/caja 50 unidades"
$ws.Range("O403").Value = "Región de O'Higgins"
$ws.Range("P403").Value = 520
$ws.Range("Q403").Value = 50
$ws.Range("R403").Value = "Hortaliza"
